# Coinranking symbol-list refresh (Fri Dec 16 10:18:57 UTC 2022, GitHub Actions).
# Updates Price (D), and for a handful of rows the Coin name (B), Link (C) and
# Volume(1h) (E) columns, on the single data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking text (e.g. "253.46"); force Text formatting
# first so Excel stores the literal string instead of silently coercing it to a
# number, then strip the formatting override back off so cell styling is untouched.
$numericTextCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D25", "D26", "D27", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49")
foreach ($addr in $numericTextCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "253.46"
$ws.Range("D3").Value = "23.66"
$ws.Range("D4").Value = "6.134"
$ws.Range("D5").Value = "0.05990"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "3.441"
$ws.Range("E6").Value = "5GateTokenGT"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "6.584"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("D8").Value = "1.321"
$ws.Range("D9").Value = "0.8014"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1515"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07921"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03350"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03095"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09294"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.585"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001680"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04792"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0006108"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "0.006228"
$ws.Range("D20").Value = "0.005713"
$ws.Range("D21").Value = "0.001079"
$ws.Range("D23").Value = "3.682"
$ws.Range("D25").Value = "0.3340"
$ws.Range("D26").Value = "0.1258"
$ws.Range("D27").Value = "0.0006502"
$ws.Range("D40").Value = "0.04461"
$ws.Range("D41").Value = "0.007053"
$ws.Range("D42").Value = "0.1072"
$ws.Range("D43").Value = "0.003367"
$ws.Range("D45").Value = "0.002471"
$ws.Range("D46").Value = "0.00005901"
$ws.Range("D48").Value = "0.7030"
$ws.Range("D49").Value = "0.1006"
$ws.Range("E49").Value = "48BOLOBOLO"

foreach ($addr in $numericTextCells) { $ws.Range($addr).ClearFormats() }
